# Update "想去人数" (interested-count) values in column F across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets to reflect
# newly scraped counts (gh-pages regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 15157
$ws1.Range("F7").Value  = 411
$ws1.Range("F9").Value  = 684
$ws1.Range("F10").Value = 15278
$ws1.Range("F12").Value = 8822
$ws1.Range("F15").Value = 73
$ws1.Range("F20").Value = 34
$ws1.Range("F24").Value = 55
$ws1.Range("F27").Value = 16
$ws1.Range("F32").Value = 37
$ws1.Range("F34").Value = 235
$ws1.Range("F38").Value = 5416

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 62

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 15157
$ws4.Range("F7").Value  = 411
$ws4.Range("F9").Value  = 684
$ws4.Range("F10").Value = 15278
$ws4.Range("F12").Value = 8822
$ws4.Range("F16").Value = 73
$ws4.Range("F21").Value = 34
$ws4.Range("F25").Value = 55
$ws4.Range("F28").Value = 16
$ws4.Range("F32").Value = 62
$ws4.Range("F35").Value = 37
$ws4.Range("F37").Value = 235
$ws4.Range("F41").Value = 5416
